$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.625493333333333
$ws.Range("H2").Value = 4.87648
$ws.Range("I2").Value = 0.1468796758507528
$ws.Range("J2").Value = 0.1468796758507528
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 52.91030366666666
$ws.Range("N2").Value = 158.730911
$ws.Range("O2").Value = 0.4161415425564564
$ws.Range("P2").Value = 0.4161415425564564
$ws.Range("Q2").Value = 86.00534587480888
$ws.Range("R2").Value = 774.0481128732799
$ws.Range("S2").Value = 0.06112273487872454
$ws.Range("T2").Value = 0.06112273487872456

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.625493333333333
$ws.Range("H3").Value = 4.87648
$ws.Range("I3").Value = 0.1468796758507528
$ws.Range("J3").Value = 0.1468796758507528
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 47.26005833333333
$ws.Range("N3").Value = 141.780175
$ws.Range("O3").Value = 0.3717021489810786
$ws.Range("P3").Value = 0.3717021489810786
$ws.Range("Q3").Value = 76.82090975377776
$ws.Range("R3").Value = 691.3881877839999
$ws.Range("S3").Value = 0.05459549115536903
$ws.Range("T3").Value = 0.05459549115536904

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.625493333333333
$ws.Range("H4").Value = 4.87648
$ws.Range("I4").Value = 0.1468796758507528
$ws.Range("J4").Value = 0.1468796758507528
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 26.97460733333333
$ws.Range("N4").Value = 80.923822
$ws.Range("O4").Value = 0.2121563084624651
$ws.Range("P4").Value = 0.2121563084624651
$ws.Range("Q4").Value = 43.84704438961778
$ws.Range("R4").Value = 394.62339950656
$ws.Range("S4").Value = 0.0311614498166592
$ws.Range("T4").Value = 0.0311614498166592

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.041193333333333
$ws.Range("H5").Value = 15.12358
$ws.Range("I5").Value = 0.4555225343081337
$ws.Range("J5").Value = 0.4555225343081337
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 52.91030366666666
$ws.Range("N5").Value = 158.730911
$ws.Range("O5").Value = 0.4161415425564564
$ws.Range("P5").Value = 0.4161415425564564
$ws.Range("Q5").Value = 266.7310701090422
$ws.Range("R5").Value = 2400.57963098138
$ws.Range("S5").Value = 0.1895618500962131
$ws.Range("T5").Value = 0.1895618500962131

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.041193333333333
$ws.Range("H6").Value = 15.12358
$ws.Range("I6").Value = 0.4555225343081337
$ws.Range("J6").Value = 0.4555225343081337
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 47.26005833333333
$ws.Range("N6").Value = 141.780175
$ws.Range("O6").Value = 0.3717021489810786
$ws.Range("P6").Value = 0.3717021489810786
$ws.Range("Q6").Value = 238.2470910029444
$ws.Range("R6").Value = 2144.2238190265
$ws.Range("S6").Value = 0.1693187049116404
$ws.Range("T6").Value = 0.1693187049116404

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.041193333333333
$ws.Range("H7").Value = 15.12358
$ws.Range("I7").Value = 0.4555225343081337
$ws.Range("J7").Value = 0.4555225343081337
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 26.97460733333333
$ws.Range("N7").Value = 80.923822
$ws.Range("O7").Value = 0.2121563084624651
$ws.Range("P7").Value = 0.2121563084624651
$ws.Range("Q7").Value = 135.9842106580844
$ws.Range("R7").Value = 1223.85789592276
$ws.Range("S7").Value = 0.09664197930028026
$ws.Range("T7").Value = 0.09664197930028028

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.400149666666667
$ws.Range("H8").Value = 13.200449
$ws.Range("I8").Value = 0.3975977898411136
$ws.Range("J8").Value = 0.3975977898411136
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 52.91030366666666
$ws.Range("N8").Value = 158.730911
$ws.Range("O8").Value = 0.4161415425564564
$ws.Range("P8").Value = 0.4161415425564564
$ws.Range("Q8").Value = 232.8132550421154
$ws.Range("R8").Value = 2095.319295379039
$ws.Range("S8").Value = 0.1654569575815188
$ws.Range("T8").Value = 0.1654569575815188

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.400149666666667
$ws.Range("H9").Value = 13.200449
$ws.Range("I9").Value = 0.3975977898411136
$ws.Range("J9").Value = 0.3975977898411136
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 47.26005833333333
$ws.Range("N9").Value = 141.780175
$ws.Range("O9").Value = 0.3717021489810786
$ws.Range("P9").Value = 0.3717021489810786
$ws.Range("Q9").Value = 207.9513299220638
$ws.Range("R9").Value = 1871.561969298575
$ws.Range("S9").Value = 0.1477879529140692
$ws.Range("T9").Value = 0.1477879529140692

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.400149666666667
$ws.Range("H10").Value = 13.200449
$ws.Range("I10").Value = 0.3975977898411136
$ws.Range("J10").Value = 0.3975977898411136
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 26.97460733333333
$ws.Range("N10").Value = 80.923822
$ws.Range("O10").Value = 0.2121563084624651
$ws.Range("P10").Value = 0.2121563084624651
$ws.Range("Q10").Value = 118.6923094662309
$ws.Range("R10").Value = 1068.230785196078
$ws.Range("S10").Value = 0.08435287934552568
$ws.Range("T10").Value = 0.08435287934552568
